$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Id = 1)
$ws.Range("B2").Value = "David"
$ws.Range("C2").Value = "Santos"
$ws.Range("D2").Value = "Lopez"
$ws.Range("E2").Value = "Perez"
$ws.Range("F2").Value = "HSFF9808KOPOIN92D2"

# Row 3 (Id = 2)
$ws.Range("B3").Value = "Armando"
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "Mora"
$ws.Range("E3").Value = "Valles"
$ws.Range("F3").Value = "MOVAAAAAA92D2"

# Row 4 (Id = 3)
$ws.Range("B4").Value = "Jorge"
$ws.Range("C4").Value = "Armando"
$ws.Range("D4").Value = "Ramirez"
$ws.Range("E4").Value = "Rocha"
$ws.Range("F4").Value = "HCH1308AAA92D2"
